# Removing less than USD 5 price from extrapolation calibration because it is just a noise
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = 115081.2387247636
$ws.Range("E6").Value = -0.01261852463264864
$ws.Range("F6").Value = 0.2106479378526508
$ws.Range("G6").Value = -0.7087787648184563
$ws.Range("H6").Value = 5.923492048363634
$ws.Range("D7").Value = 115729.4014244649
$ws.Range("E7").Value = -0.01549382194285488
$ws.Range("F7").Value = 0.210189590555609
$ws.Range("G7").Value = -0.5168589211625243
$ws.Range("H7").Value = 5.223422583104654
$ws.Range("D8").Value = 115902.4441511091
$ws.Range("E8").Value = -0.03242949234361994
$ws.Range("F8").Value = 0.2085967379298454
$ws.Range("G8").Value = -1.348204895146933
$ws.Range("H8").Value = 8.870920615642103
$ws.Range("D9").Value = 117402.6028498795
$ws.Range("E9").Value = -0.0636740399468614
$ws.Range("F9").Value = 0.3356250937238441
$ws.Range("G9").Value = -1.827256894593787
$ws.Range("H9").Value = 10.98182625486903
$ws.Range("D10").Value = 118988.3896664457
$ws.Range("E10").Value = -0.09788690023333591
$ws.Range("F10").Value = 0.4150687660883884
$ws.Range("G10").Value = -1.879437335611029
$ws.Range("H10").Value = 9.608403840296544
$ws.Range("D11").Value = 120880.7258083007
$ws.Range("E11").Value = -0.1667754622575358
$ws.Range("F11").Value = 0.6951836192414919
$ws.Range("G11").Value = -2.471494357357156
$ws.Range("H11").Value = 11.58301120960439
$ws.Range("D15").Value = 111922.032564998
$ws.Range("E15").Value = 0.03152131171749281
$ws.Range("F15").Value = 0.1428945784418704
$ws.Range("G15").Value = -0.7851727356124144
$ws.Range("H15").Value = 6.075221095169978
$ws.Range("D17").Value = 111895.3534178812
$ws.Range("E17").Value = 0.02121193946449751
$ws.Range("F17").Value = 0.1562058505465646
$ws.Range("G17").Value = -0.7395612067214531
$ws.Range("H17").Value = 5.425419922712988

Write-Host 'Applied calibration updates (removed less than USD 5 price noise)'
